$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("location-page.html")

# New translation rows: origin key, English, Ukrainian, Russian
$ws.Cells.Item(13, 1).Value = "PHX_REFRESH_LIST"
$ws.Cells.Item(13, 2).Value = "Refresh the list"
$ws.Cells.Item(13, 3).Value = "Оновити список"
$ws.Cells.Item(13, 4).Value = "Обновить список"

$ws.Cells.Item(14, 1).Value = "PHX_START_LIST"
$ws.Cells.Item(14, 2).Value = "Start list"
$ws.Cells.Item(14, 3).Value = "Старт списку"
$ws.Cells.Item(14, 4).Value = "Старт списка"

$ws.Cells.Item(15, 1).Value = "PHX_DELETE_ARENA"
$ws.Cells.Item(15, 2).Value = "Delete arena"
$ws.Cells.Item(15, 3).Value = "Видалити арена"
$ws.Cells.Item(15, 4).Value = "Удалить арена"
